# #44 equip effect update
# The "护甲" (Armor) sheet's VitP (column I) base effect value is bumped
# from 0 to 100 for every equipment row (rows 4-32).

$wb = $excel.ActiveWorkbook

# "护甲" (Armor) is the 3rd tab (头盔, 武器, 护甲, 饰品).
$wsArmor = $wb.Worksheets.Item(3)
$wsArmor.Range("I4:I32").Value = 100

# Leave the cursor where the editing session ended on the Armor sheet.
$wsArmor.Range("I6").Select()

# Finish back on the "头盔" (Helmet) sheet, matching the workbook's last
# active tab/selection when the file was saved.
$wsHelmet = $wb.Worksheets.Item(1)
$wsHelmet.Range("H8").Select()
